$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "image"
$ws.Range("F2").Value = "C:\LZ\natural\natural_119.jpg"

$ws.Range("F2").Select()
